# "Refined metadata to be additional tab": the workbook gains a second
# worksheet ("metadata") carrying the PanelApp query metadata that used to
# live only outside the sheet, and the "data" sheet's per-row time_taken
# column (F) is refreshed to the timestamps recorded by the new query run.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Refresh "data" sheet's time_taken column (F2:F71) -----------------
$timestamps = @(
    "2021-10-05 14:21:02.530665",
    "2021-10-05 14:21:02.530673",
    "2021-10-05 14:21:02.530676",
    "2021-10-05 14:21:02.530679",
    "2021-10-05 14:21:02.530682",
    "2021-10-05 14:21:02.530685",
    "2021-10-05 14:21:02.530688",
    "2021-10-05 14:21:02.530690",
    "2021-10-05 14:21:02.530693",
    "2021-10-05 14:21:02.530696",
    "2021-10-05 14:21:02.530698",
    "2021-10-05 14:21:02.530701",
    "2021-10-05 14:21:02.530704",
    "2021-10-05 14:21:02.530706",
    "2021-10-05 14:21:02.530709",
    "2021-10-05 14:21:02.530712",
    "2021-10-05 14:21:02.530715",
    "2021-10-05 14:21:02.530717",
    "2021-10-05 14:21:02.530720",
    "2021-10-05 14:21:02.530722",
    "2021-10-05 14:21:02.530725",
    "2021-10-05 14:21:02.530728",
    "2021-10-05 14:21:02.530730",
    "2021-10-05 14:21:02.530733",
    "2021-10-05 14:21:02.530736",
    "2021-10-05 14:21:02.530738",
    "2021-10-05 14:21:02.530741",
    "2021-10-05 14:21:02.530744",
    "2021-10-05 14:21:02.530747",
    "2021-10-05 14:21:02.530749",
    "2021-10-05 14:21:02.530752",
    "2021-10-05 14:21:02.530754",
    "2021-10-05 14:21:02.530757",
    "2021-10-05 14:21:02.530760",
    "2021-10-05 14:21:02.530763",
    "2021-10-05 14:21:02.530765",
    "2021-10-05 14:21:02.530768",
    "2021-10-05 14:21:02.530771",
    "2021-10-05 14:21:02.530774",
    "2021-10-05 14:21:02.530776",
    "2021-10-05 14:21:02.530780",
    "2021-10-05 14:21:02.530782",
    "2021-10-05 14:21:02.530785",
    "2021-10-05 14:21:02.530788",
    "2021-10-05 14:21:02.530790",
    "2021-10-05 14:21:02.530793",
    "2021-10-05 14:21:02.530796",
    "2021-10-05 14:21:02.530798",
    "2021-10-05 14:21:02.530801",
    "2021-10-05 14:21:02.530804",
    "2021-10-05 14:21:02.530806",
    "2021-10-05 14:21:02.530809",
    "2021-10-05 14:21:02.530812",
    "2021-10-05 14:21:02.530814",
    "2021-10-05 14:21:02.530817",
    "2021-10-05 14:21:02.530820",
    "2021-10-05 14:21:02.530823",
    "2021-10-05 14:21:02.530825",
    "2021-10-05 14:21:02.530828",
    "2021-10-05 14:21:02.530831",
    "2021-10-05 14:21:02.530833",
    "2021-10-05 14:21:02.530836",
    "2021-10-05 14:21:02.530839",
    "2021-10-05 14:21:02.530841",
    "2021-10-05 14:21:02.530846",
    "2021-10-05 14:21:02.530849",
    "2021-10-05 14:21:02.530852",
    "2021-10-05 14:21:02.530855",
    "2021-10-05 14:21:02.530857",
    "2021-10-05 14:21:02.530860"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the new "metadata" worksheet, placed right after "data" -------
$metadata = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$metadata.Name = "metadata"

# Header row
$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"
$metadata.Range("G1").Value = "panel_get_request"

# Reuse the "data" sheet's existing bold/border/center header style
# instead of building a new font/style combination.
$data.Range("B1:F1").Copy()
$metadata.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$data.Range("B1").Copy()
$metadata.Range("G1").PasteSpecial(-4122)

# Data row
$metadata.Range("A2").Value = 0
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)  # reuse same style as column A in "data"

$metadata.Range("B2").Value = "Ichthyosis and erythrokeratoderma"
$metadata.Range("C2").Value = 555
# data_version is stored as text ("1.68"), not a number - force text format first
$metadata.Range("D2").NumberFormat = "@"
$metadata.Range("D2").Value = "1.68"
$metadata.Range("E2").Value = "2021-09-30T14:21:23.768335Z"
$metadata.Range("F2").Value = "2021-10-05 14:21:02.527439"
$metadata.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/555/?format=json"

# Keep "data" as the active/selected sheet, matching the original workbook view.
$data.Activate()
